$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "steps" (sheet3): lowercase existing step titles and append new steps
# ---------------------------------------------------------------------------
$steps = $wb.Worksheets.Item("steps")

$steps.Range("A2").Value = "add"
$steps.Range("A3").Value = "let boil"
$steps.Range("A4").Value = "let simmer"
$steps.Range("A5").Value = "let cook"
$steps.Range("A6").Value = "boil"
$steps.Range("A7").Value = "mix"

$steps.Range("A8").Value = "rinse"
$steps.Range("C8").Value = "wait"

$steps.Range("A9").Value = "fluff"
$steps.Range("B9").Value = "Fluff with fork"
$steps.Range("C9").Value = "wait"

$steps.Range("A10").Value = "serve"
$steps.Range("C10").Value = "wait"

# ---------------------------------------------------------------------------
# Sheet "recipe_steps" (sheet6): add an "order" column and populate the
# steps for the two recipes (white basmati rice / long grain brown rice)
# ---------------------------------------------------------------------------
$rs = $wb.Worksheets.Item("recipe_steps")

# Header row - insert "order" before "seconds", and add "weight"
$rs.Range("D1").Value = "order"
$rs.Range("E1").Value = "seconds"
$rs.Range("F1").Value = "weight"

# --- white basmati rice ---
$rs.Range("A2").Value = "basmati rice"
$rs.Range("B2").Value = "white basmati rice"
$rs.Range("C2").Value = "add"
$rs.Range("D2").Value = 1
$rs.Range("F2").Value = 200

$rs.Range("A3").Value = "water"
$rs.Range("B3").Value = "white basmati rice"
$rs.Range("C3").Value = "add"
$rs.Range("D3").Value = 2
$rs.Range("F3").Value = 295

$rs.Range("A4").Value = "table salt"
$rs.Range("B4").Value = "white basmati rice"
$rs.Range("C4").Value = "add"
$rs.Range("D4").Value = 3
$rs.Range("F4").Value = 3

$rs.Range("B5").Value = "white basmati rice"
$rs.Range("C5").Value = "boil"
$rs.Range("D5").Value = 4

$rs.Range("B6").Value = "white basmati rice"
$rs.Range("C6").Value = "let boil"
$rs.Range("D6").Value = 5
$rs.Range("E6").Value = 1200

$rs.Range("B7").Value = "white basmati rice"
$rs.Range("C7").Value = "fluff"
$rs.Range("D7").Value = 6

$rs.Range("B8").Value = "white basmati rice"
$rs.Range("C8").Value = "serve"
$rs.Range("D8").Value = 7

# --- long grain brown rice ---
$rs.Range("A9").Value = "long grain brown rice"
$rs.Range("B9").Value = "long grain brown rice"
$rs.Range("C9").Value = "add"
$rs.Range("D9").Value = 1
$rs.Range("F9").Value = 190

$rs.Range("A10").Value = "water"
$rs.Range("B10").Value = "long grain brown rice"
$rs.Range("C10").Value = "add"
$rs.Range("D10").Value = 2
$rs.Range("F10").Value = 475

$rs.Range("A11").Value = "table salt"
$rs.Range("B11").Value = "long grain brown rice"
$rs.Range("C11").Value = "add"
$rs.Range("D11").Value = 3
$rs.Range("F11").Value = 5.5

$rs.Range("B12").Value = "long grain brown rice"
$rs.Range("C12").Value = "boil"
$rs.Range("D12").Value = 4
$rs.Range("E12").Value = 2700

$rs.Range("B13").Value = "long grain brown rice"
$rs.Range("C13").Value = "let boil"
$rs.Range("D13").Value = 5

$rs.Range("B14").Value = "long grain brown rice"
$rs.Range("C14").Value = "serve"
$rs.Range("D14").Value = 6

# ---------------------------------------------------------------------------
# Active tab / active sheet -> "steps" (index 2, 0-based) ends up active
# ---------------------------------------------------------------------------
$steps.Activate()

"done"
